$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the two Spine nodes (103/104 -> 201/202)
$ws.Range("E4").Value = "Spine-201"
$ws.Range("F4").Value = 201
$ws.Range("E5").Value = "Spine-202"
$ws.Range("F5").Value = 202

# Fix column B (Node Type/Role) data validation: the first rule on B2:B5
# had a stale list; align it with the second rule's allowed values.
$ws.Range("B2:B5").Validation.Formula1 = "spine,leaf,unspecified"
